$wb = $excel.ActiveWorkbook

# --- Update sigma_010 (sheet2) ---
$ws2 = $wb.Worksheets.Item("sigma_010")
$ws2.Cells.Item(2, 2).Value = 27.85093759301251
$ws2.Cells.Item(2, 3).Value = 31.30499278530378
$ws2.Cells.Item(3, 2).Value = 27.81658219922084
$ws2.Cells.Item(3, 3).Value = 31.26221328012794
$ws2.Cells.Item(4, 2).Value = 27.7938933766297
$ws2.Cells.Item(4, 3).Value = 31.27200857524841
$ws2.Cells.Item(5, 2).Value = 27.8041697442676
$ws2.Cells.Item(5, 3).Value = 31.28438020501203
$ws2.Cells.Item(6, 2).Value = 27.80005135302297
$ws2.Cells.Item(6, 3).Value = 31.29535990136297
$ws2.Cells.Item(7, 2).Value = 27.82345045813658
$ws2.Cells.Item(7, 3).Value = 31.26162660740683
$ws2.Cells.Item(8, 2).Value = 27.85300797181988
$ws2.Cells.Item(8, 3).Value = 31.26090805584424
$ws2.Cells.Item(9, 2).Value = 27.79091590194805
$ws2.Cells.Item(9, 3).Value = 31.27829912349581
$ws2.Cells.Item(10, 2).Value = 27.83600223158515
$ws2.Cells.Item(10, 3).Value = 31.28862073911401
$ws2.Cells.Item(11, 2).Value = 27.8228161759004
$ws2.Cells.Item(11, 3).Value = 31.29575395987793
$ws2.Cells.Item(12, 2).Value = 27.81918270055437
$ws2.Cells.Item(12, 3).Value = 31.2804163232794

# --- Update sigma_025 (sheet3) ---
$ws3 = $wb.Worksheets.Item("sigma_025")
$ws3.Cells.Item(2, 2).Value = 19.70500473168443
$ws3.Cells.Item(2, 3).Value = 27.76008228379396
$ws3.Cells.Item(3, 2).Value = 19.69010847023509
$ws3.Cells.Item(3, 3).Value = 27.75731857934534
$ws3.Cells.Item(4, 2).Value = 19.68805822976169
$ws3.Cells.Item(4, 3).Value = 27.74161545163189
$ws3.Cells.Item(5, 2).Value = 19.6919575949737
$ws3.Cells.Item(5, 3).Value = 27.69929462883706
$ws3.Cells.Item(6, 2).Value = 19.68342974970465
$ws3.Cells.Item(6, 3).Value = 27.76190811334769
$ws3.Cells.Item(7, 2).Value = 19.68750592871477
$ws3.Cells.Item(7, 3).Value = 27.78329529707954
$ws3.Cells.Item(8, 2).Value = 19.69425374097414
$ws3.Cells.Item(8, 3).Value = 27.77822413744289
$ws3.Cells.Item(9, 2).Value = 19.67751416913032
$ws3.Cells.Item(9, 3).Value = 27.74053838164122
$ws3.Cells.Item(10, 2).Value = 19.67754537564062
$ws3.Cells.Item(10, 3).Value = 27.77332333341515
$ws3.Cells.Item(11, 2).Value = 19.67284770738014
$ws3.Cells.Item(11, 3).Value = 27.74028356326328
$ws3.Cells.Item(12, 2).Value = 19.68682256981996
$ws3.Cells.Item(12, 3).Value = 27.7535883769798

# --- Add new sheet sigma_050 (sheet4) at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "sigma_050"

$ws4.Cells.Item(1, 1).Value = "Rows"
$ws4.Cells.Item(1, 2).Value = "Noisy"
$ws4.Cells.Item(1, 3).Value = "NLM-LBP"
$ws4.Cells.Item(2, 1).Value = 0
$ws4.Cells.Item(2, 2).Value = 14.67218905703577
$ws4.Cells.Item(2, 3).Value = 22.5469566443962
$ws4.Cells.Item(3, 1).Value = 1
$ws4.Cells.Item(3, 2).Value = 14.68277900071762
$ws4.Cells.Item(3, 3).Value = 22.56495161203066
$ws4.Cells.Item(4, 1).Value = 2
$ws4.Cells.Item(4, 2).Value = 14.69089101910216
$ws4.Cells.Item(4, 3).Value = 22.52866168101401
$ws4.Cells.Item(5, 1).Value = 3
$ws4.Cells.Item(5, 2).Value = 14.69064886449247
$ws4.Cells.Item(5, 3).Value = 22.53348694456143
$ws4.Cells.Item(6, 1).Value = 4
$ws4.Cells.Item(6, 2).Value = 14.69355985260382
$ws4.Cells.Item(6, 3).Value = 22.56257848337445
$ws4.Cells.Item(7, 1).Value = 5
$ws4.Cells.Item(7, 2).Value = 14.69414521400686
$ws4.Cells.Item(7, 3).Value = 22.52944917738828
$ws4.Cells.Item(8, 1).Value = 6
$ws4.Cells.Item(8, 2).Value = 14.67408945050213
$ws4.Cells.Item(8, 3).Value = 22.54715835203611
$ws4.Cells.Item(9, 1).Value = 7
$ws4.Cells.Item(9, 2).Value = 14.69549093456045
$ws4.Cells.Item(9, 3).Value = 22.5189632524404
$ws4.Cells.Item(10, 1).Value = 8
$ws4.Cells.Item(10, 2).Value = 14.67949299936367
$ws4.Cells.Item(10, 3).Value = 22.55797782978884
$ws4.Cells.Item(11, 1).Value = 9
$ws4.Cells.Item(11, 2).Value = 14.68606969036923
$ws4.Cells.Item(11, 3).Value = 22.58686847416819
$ws4.Cells.Item(12, 1).Value = "Média"
$ws4.Cells.Item(12, 2).Value = 14.68593560827542
$ws4.Cells.Item(12, 3).Value = 22.54770524511986

Write-Host "Done."